# NIT-9005157676.xlsx - "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
#
# Semantic edits captured by this change:
#  1. The "Periodo Mora" value shown for the last three detail rows (17-19)
#     changes from 2508 to 2509.
#  2. The "Periodo Mora" column (E16:E19) gets horizontal-center alignment
#     applied (it previously used the sheet's default/general alignment).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the period value (2508 -> 2509) on the detail rows that carry it.
$ws.Range("E17").Value = "2509"
$ws.Range("E18").Value = "2509"
$ws.Range("E19").Value = "2509"

# 2) Center-align the whole "Periodo Mora" column for the detail rows.
$ws.Range("E16:E19").HorizontalAlignment = -4108
